$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, matching the style used by the other header cells (B1:G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add a "Save" value of 0 for each data row (2-12), matching column G's (unstyled) cells
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 8).Value = 0
}
